# The sheet's kitchens-categorization trial list was regenerated: one
# target trial (stimuli/img_r77yy.png) was dropped, the remaining 26
# trials were reshuffled into a new presentation order, and trial_block
# (E) / trial_total (F) were renumbered sequentially to match. All other
# per-stimulus fields are unchanged from their original values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    (16, 'categorization', 5, 2, 1, 217, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_ye5sl.png', 53.2258064516129, 34.45161290322581, 43.83870967741936, 31, 2, 2, 2),
    (16, 'categorization', 5, 2, 2, 218, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_9mky8.png', 84.32352941176471, 65.17647058823529, 74.75, 34, 9, 9, 9),
    (16, 'categorization', 5, 2, 3, 219, 'kitchens', 'living_rooms', 'distractor', $null, 'f', 'stimuli/img_jpjeg.png', 90.90697674418605, 74.3953488372093, 82.65116279069767, 43, 10, 10, 10),
    (16, 'categorization', 5, 2, 4, 220, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_uwv6y.png', 78.88888888888889, 59.30555555555556, 69.09722222222223, 36, 8, 8, 8),
    (16, 'categorization', 5, 2, 5, 221, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_nyv2b.png', 11.91176470588235, 6.852941176470588, 9.382352941176471, 34, 1, 1, 1),
    (16, 'categorization', 5, 2, 6, 222, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_7wul8.png', 43.03030303030303, 25.54545454545455, 34.28787878787879, 33, 1, 1, 1),
    (16, 'categorization', 5, 2, 7, 223, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_d8xbu.png', 91.36363636363636, 73.18181818181819, 82.27272727272728, 33, 10, 10, 10),
    (16, 'categorization', 5, 2, 8, 224, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_t90e2.png', 83.0625, 61.96875, 72.515625, 32, 9, 9, 9),
    (16, 'categorization', 5, 2, 9, 225, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_inqod.png', 70.84848484848484, 50.63636363636363, 60.74242424242424, 33, 5, 5, 5),
    (16, 'categorization', 5, 2, 10, 226, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_p3hpc.png', 72.83333333333333, 52.22222222222222, 62.52777777777777, 36, 6, 6, 6),
    (16, 'categorization', 5, 2, 11, 227, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_463mq.png', 51.35294117647059, 30.20588235294118, 40.77941176470588, 34, 2, 2, 2),
    (16, 'categorization', 5, 2, 12, 228, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_wyl6z.png', 59.8235294117647, 36.23529411764706, 48.02941176470588, 34, 3, 3, 3),
    (16, 'categorization', 5, 2, 13, 229, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_cnyac.png', 69.14705882352941, 47.8235294117647, 58.48529411764706, 34, 5, 5, 5),
    (16, 'categorization', 5, 2, 14, 230, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_yeh72.png', 68.66666666666667, 45.21212121212121, 56.93939393939394, 33, 4, 4, 4),
    (16, 'categorization', 5, 2, 15, 231, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_3gm8h.png', 65.07894736842105, 43.92105263157895, 54.5, 38, 4, 4, 4),
    (16, 'categorization', 5, 2, 16, 232, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_60242.png', 78.33333333333333, 57.57575757575758, 67.95454545454545, 33, 7, 7, 7),
    (16, 'categorization', 5, 2, 17, 233, 'kitchens', 'living_rooms', 'distractor', $null, 'f', 'stimuli/img_53nbn.png', 73.28888888888889, 51.15555555555556, 62.22222222222223, 45, 6, 6, 6),
    (16, 'categorization', 5, 2, 18, 234, 'kitchens', 'living_rooms', 'distractor', $null, 'f', 'stimuli/img_x9w7o.png', 92.38888888888889, 72.94444444444444, 82.66666666666666, 36, 10, 10, 10),
    (16, 'categorization', 5, 2, 19, 235, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_a8wvq.png', 86.25925925925925, 66.25925925925925, 76.25925925925925, 27, 10, 10, 10),
    (16, 'categorization', 5, 2, 20, 236, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_aplao.png', 64.09090909090909, 40.75757575757576, 52.42424242424242, 33, 3, 3, 3),
    (16, 'categorization', 5, 2, 21, 237, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_iyxnj.png', 75.30555555555556, 54.33333333333334, 64.81944444444444, 36, 6, 6, 6),
    (16, 'categorization', 5, 2, 22, 238, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_eatdk.png', 81.40625, 61.375, 71.390625, 32, 8, 8, 8),
    (16, 'categorization', 5, 2, 23, 239, 'kitchens', 'bedrooms', 'distractor', $null, 'f', 'stimuli/img_scrdm.png', 78.675, 57.9, 68.28749999999999, 40, 7, 7, 7),
    (16, 'categorization', 5, 2, 24, 240, 'kitchens', 'kitchens', 'target', $null, 'j', 'stimuli/img_ce9vx.png', 75.90909090909091, 57.12121212121212, 66.51515151515152, 33, 7, 7, 7),
    (16, 'categorization', 5, 2, 25, 241, 'kitchens', 'living_rooms', 'distractor', $null, 'f', 'stimuli/img_g13d5.png', 73, 51.51111111111111, 62.25555555555556, 45, 6, 6, 6),
    (16, 'categorization', 5, 2, 26, 242, 'kitchens', 'bedrooms', 'distractor', $null, 'f', 'stimuli/img_kn0we.png', 80.15909090909091, 56.68181818181818, 68.42045454545455, 44, 7, 7, 7)
)

$numRows = $data.Length
$numCols = 19

$arr = New-Object 'object[,]' $numRows, $numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $arr[$r, $c] = $data[$r][$c]
    }
}

# Overwrite the 26 data rows (A2:S27) with the regenerated trial list...
$ws.Range("A2:S27").Value = $arr

# ...and drop the now-removed 28th row so the sheet shrinks from
# A1:S28 to A1:S27.
$ws.Range("A28:S28").ClearContents()
